# AMEK Tidsredovisning - add two new time-tracking rows below the existing
# "Ritningsnotiser" note block (rows 45-46), matching the sheet's existing
# layout/styles for that block (date, hours, purpose).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: 2020-03-31, 2h, "Ändra, utskrift, sortering , visa av ritningsnotiser"
$ws.Range("A47").Value = 43921
$ws.Range("B47").Value = 2
$ws.Range("C47").Value = "Ändra, utskrift, sortering , visa av ritningsnotiser"

# Row 48: 2020-04-04, 0.5h, "Extra visning av ritningsnotering i Orderfönstret"
$ws.Range("A48").Value = 43925
$ws.Range("B48").Value = 0.5
$ws.Range("C48").Value = "Extra visning av ritningsnotering i Orderfönstret"

# Match the author's final selection/active cell after data entry
$ws.Range("B48").Select()
